$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.780.74"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "1.631.92"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.72%  "
$c = $ws.Range("D5")
$c.Value = "'211.69"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("E6").Value = "  +0.03%  "
$c = $ws.Range("D7")
$c.Value = "'0.994"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.80%  "
$c = $ws.Range("D8")
$c.Value = "'23.35"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("E9").Value = "  -1.98%  "
$c = $ws.Range("D10")
$c.Value = "'0.0611"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.03%  "
$c = $ws.Range("D11")
$c.Value = "'0.0879"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("D12").Value = "1.863.02"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "1.632.66"
$ws.Range("E13").Value = "  +0.27%  "
$c = $ws.Range("D14")
$c.Value = "'4.06"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("E15").Value = "  +2.10%  "
$c = $ws.Range("D16")
$c.Value = "'65.29"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").Value = "27.775.61"
$ws.Range("E17").Value = "  +1.20%  "
$c = $ws.Range("D18")
$c.Value = "'233.14"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.95%  "
$ws.Range("D19").Value = "0.0₃0719"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  -0.80%  "
$c = $ws.Range("D22")
$c.Value = "'10.45"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.02%  "
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("E24").Value = "  -2.67%  "
$c = $ws.Range("D25")
$c.Value = "'151.64"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.43%  "
$c = $ws.Range("D26")
$c.Value = "'6.86"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.03%  "
$c = $ws.Range("D27")
$c.Value = "'15.60"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  -0.25%  "
$c = $ws.Range("D29")
$c.Value = "'0.995"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.75%  "
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +2.26%  "
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("D34").Value = "1.402.94"
$ws.Range("E34").Value = "  -4.22%  "
$ws.Range("E35").Value = "  +2.13%  "
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("E38").Value = "  +0.18%  "
$ws.Range("E39").Value = "  -0.48%  "
$c = $ws.Range("D40")
$c.Value = "'0.914"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  -0.88%  "
$c = $ws.Range("D43")
$c.Value = "'66.90"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.30%  "
$c = $ws.Range("D44")
$c.Value = "'1.85"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +6.51%  "
$ws.Range("E45").Value = "  +1.95%  "
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "1.772.65"
$ws.Range("E47").Value = "  +0.36%  "
$c = $ws.Range("D48")
$c.Value = "'87.41"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'0.0996"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "
$c = $ws.Range("D50")
$c.Value = "'0.0505"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.27%  "
$ws.Range("E51").Value = "  -1.47%  "
